$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "time"
$ws.Range("C1").Value = "totalTime"

# Data rows
$ws.Range("A2").Value = "chair"
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 10

$ws.Range("A3").Value = "table"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 10

$ws.Range("A4").Value = "tv"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 10

$ws.Range("A5").Value = "microwave"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 10

$ws.Range("A6").Value = "sofa"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 10
